$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 116.9511416666667
$ws.Range("H2").Value = 350.853425
$ws.Range("I2").Value = 0.411881549221027
$ws.Range("J2").Value = 0.411881549221027
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.59487733333333
$ws.Range("N2").Value = 67.784632
$ws.Range("O2").Value = 0.7395019553569895
$ws.Range("P2").Value = 0.7395019553569895
$ws.Range("Q2").Value = 2642.496699951622
$ws.Range("R2").Value = 23782.4702995646
$ws.Range("S2").Value = 0.3045872110244155
$ws.Range("T2").Value = 0.3045872110244155
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 116.9511416666667
$ws.Range("H3").Value = 350.853425
$ws.Range("I3").Value = 0.411881549221027
$ws.Range("J3").Value = 0.411881549221027
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.3045986666666667
$ws.Range("N3").Value = 0.9137960000000001
$ws.Range("O3").Value = 0.00996913177602551
$ws.Range("P3").Value = 0.00996913177602551
$ws.Range("Q3").Value = 35.62316181681111
$ws.Range("R3").Value = 320.6084563513
$ws.Range("S3").Value = 0.004106101440297955
$ws.Range("T3").Value = 0.004106101440297955
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 116.9511416666667
$ws.Range("H4").Value = 350.853425
$ws.Range("I4").Value = 0.411881549221027
$ws.Range("J4").Value = 0.411881549221027
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.654706
$ws.Range("N4").Value = 22.964118
$ws.Range("O4").Value = 0.2505289128669849
$ws.Range("P4").Value = 0.2505289128669849
$ws.Range("Q4").Value = 895.2266058226834
$ws.Range("R4").Value = 8057.03945240415
$ws.Range("S4").Value = 0.1031882367563134
$ws.Range("T4").Value = 0.1031882367563134
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 133.0753813333333
$ws.Range("H5").Value = 399.226144
$ws.Range("I5").Value = 0.468668312644395
$ws.Range("J5").Value = 0.468668312644395
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.59487733333333
$ws.Range("N5").Value = 67.784632
$ws.Range("O5").Value = 0.7395019553569895
$ws.Range("P5").Value = 0.7395019553569895
$ws.Range("Q5").Value = 3006.821917313223
$ws.Range("R5").Value = 27061.39725581901
$ws.Range("S5").Value = 0.346581133614391
$ws.Range("T5").Value = 0.346581133614391
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 133.0753813333333
$ws.Range("H6").Value = 399.226144
$ws.Range("I6").Value = 0.468668312644395
$ws.Range("J6").Value = 0.468668312644395
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.3045986666666667
$ws.Range("N6").Value = 0.9137960000000001
$ws.Range("O6").Value = 0.00996913177602551
$ws.Range("P6").Value = 0.00996913177602551
$ws.Range("Q6").Value = 40.53458372029156
$ws.Range("R6").Value = 364.811253482624
$ws.Range("S6").Value = 0.004672216167999497
$ws.Range("T6").Value = 0.004672216167999497
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 133.0753813333333
$ws.Range("H7").Value = 399.226144
$ws.Range("I7").Value = 0.468668312644395
$ws.Range("J7").Value = 0.468668312644395
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.654706
$ws.Range("N7").Value = 22.964118
$ws.Range("O7").Value = 0.2505289128669849
$ws.Range("P7").Value = 0.2505289128669849
$ws.Range("Q7").Value = 1018.652919944555
$ws.Range("R7").Value = 9167.876279500992
$ws.Range("S7").Value = 0.1174149628620045
$ws.Range("T7").Value = 0.1174149628620045
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 33.91710566666666
$ws.Range("H8").Value = 101.751317
$ws.Range("I8").Value = 0.119450138134578
$ws.Range("J8").Value = 0.119450138134578
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 22.59487733333333
$ws.Range("N8").Value = 67.784632
$ws.Range("O8").Value = 0.7395019553569895
$ws.Range("P8").Value = 0.7395019553569895
$ws.Range("Q8").Value = 766.3528420400381
$ws.Range("R8").Value = 6897.175578360344
$ws.Range("S8").Value = 0.08833361071818295
$ws.Range("T8").Value = 0.08833361071818295
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 33.91710566666666
$ws.Range("H9").Value = 101.751317
$ws.Range("I9").Value = 0.119450138134578
$ws.Range("J9").Value = 0.119450138134578
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.3045986666666667
$ws.Range("N9").Value = 0.9137960000000001
$ws.Range("O9").Value = 0.00996913177602551
$ws.Range("P9").Value = 0.00996913177602551
$ws.Range("Q9").Value = 10.33110516325911
$ws.Range("R9").Value = 92.97994646933201
$ws.Range("S9").Value = 0.001190814167728059
$ws.Range("T9").Value = 0.001190814167728059
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 33.91710566666666
$ws.Range("H10").Value = 101.751317
$ws.Range("I10").Value = 0.119450138134578
$ws.Range("J10").Value = 0.119450138134578
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.654706
$ws.Range("N10").Value = 22.964118
$ws.Range("O10").Value = 0.2505289128669849
$ws.Range("P10").Value = 0.2505289128669849
$ws.Range("Q10").Value = 259.6254722492673
$ws.Range("R10").Value = 2336.629250243406
$ws.Range("S10").Value = 0.02992571324866701
$ws.Range("T10").Value = 0.02992571324866701
